$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.287.82"
$ws.Range("E2").Value = "  +0.24%  "

$ws.Range("D3").Value = "1.870.79"
$ws.Range("E3").Value = "  +0.19%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7082"
$ws.Range("E5").Value = "  -0.56%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "241.66"
$ws.Range("E6").Value = "  +0.00%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9999"
$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07788"
$ws.Range("E8").Value = "  +1.42%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3092"
$ws.Range("E9").Value = "  -0.85%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "25.00"
$ws.Range("E10").Value = "  +1.08%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08394"
$ws.Range("E11").Value = "  +0.26%  "

$ws.Range("D12").Value = "1.876.00"
$ws.Range("E12").Value = "  +0.05%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.241"
$ws.Range("E13").Value = "  +0.25%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7103"
$ws.Range("E14").Value = "  -0.21%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "90.96"
$ws.Range("E15").Value = "  -0.50%  "

$ws.Range("D16").Value = "29.295.56"
$ws.Range("E16").Value = "  +0.17%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.071"
$ws.Range("E17").Value = "  +1.91%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008172"
$ws.Range("E18").Value = "  +4.28%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "239.57"
$ws.Range("E19").Value = "  -1.91%  "

$ws.Range("E20").Value = "  +0.83%  "

$ws.Range("D21").Value = "2.111.26"
$ws.Range("E21").Value = "  -0.17%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9998"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.749"
$ws.Range("E23").Value = "  -1.57%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9998"
$ws.Range("E24").Value = "  -0.02%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1592"
$ws.Range("E25").Value = "  -2.36%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.10"
$ws.Range("E26").Value = "  -0.03%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.005"
$ws.Range("E27").Value = "  +0.56%  "

$ws.Range("E28").Value = "  -0.39%  "

$ws.Range("E29").Value = "  +0.20%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.390"
$ws.Range("E30").Value = "  -0.38%  "

$ws.Range("E31").Value = "  -1.52%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.286"
$ws.Range("E32").Value = "  +0.67%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05334"
$ws.Range("E33").Value = "  +3.52%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.938"
$ws.Range("E34").Value = "  +1.25%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.176"
$ws.Range("E35").Value = "  +0.68%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7445"
$ws.Range("E36").Value = "  -6.23%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.699"
$ws.Range("E37").Value = "  +0.46%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01867"
$ws.Range("E38").Value = "  +0.51%  "

$ws.Range("D39").Value = "1.227.82"
$ws.Range("E39").Value = "  +5.78%  "

$ws.Range("E40").Value = "  +0.46%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.556"
$ws.Range("E41").Value = "  +3.79%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8842"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "109.03"
$ws.Range("E43").Value = "  +5.58%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "72.33"
$ws.Range("E44").Value = "  -1.77%  "

$ws.Range("E45").Value = "  +0.02%  "

$ws.Range("D46").Value = "2.013.86"
$ws.Range("E46").Value = "  +0.08%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5191"
$ws.Range("E47").Value = "  -0.07%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.788"
$ws.Range("E48").Value = "  +0.44%  "

$ws.Range("E49").Value = "  +2.04%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.411"
$ws.Range("E50").Value = "  +0.64%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4310"
$ws.Range("E51").Value = "  +0.26%  "

